{"js": "// The \"Terreno de juego\" bullet originally reads:\n//   \"...encajaba mas con la tem\u00e1tica de Mario Bross con la que se quer\u00eda...\"\n// The edit inserts \"los juegos de \" right before \"Mario Bross\" so it reads:\n//   \"...encajaba mas con la tem\u00e1tica de los juegos de Mario Bross con la que...\"\nconst body = context.document.body;\n\n// Anchor the search on a long, unique phrase so we don't accidentally\n// match the (identically worded) sentence fragments that appear elsewhere\n// in the document.\nconst results = body.search(\"Mario Bross con la que se quer\u00eda dotar al juego\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the target sentence containing \"Mario Bross\".');\n}\n\n// Insert the new text immediately before the matched range (i.e. right\n// before \"Mario Bross\"), leaving everything else in the sentence intact.\nconst target = results.items[0];\ntarget.insertText(\"los juegos de \", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# The \"Terreno de juego\" bullet originally reads:\n#   \"...encajaba mas con la tem\u00e1tica de Mario Bross con la que se quer\u00eda...\"\n# The edit inserts \"los juegos de \" right before \"Mario Bross\" so it reads:\n#   \"...encajaba mas con la tem\u00e1tica de los juegos de Mario Bross con la que...\"\n$d = $word.ActiveDocument\n\n# Anchor the search on a long, unique phrase so we don't accidentally match\n# the (identically worded) sentence fragments that appear elsewhere in the\n# document.\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"Mario Bross con la que se quer\u00eda dotar al juego\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find the target sentence containing \"Mario Bross\".'\n}\n\n# Collapse the found range to its start (right before \"Mario Bross\") and\n# insert the new text there, leaving the rest of the sentence intact.\n$insertRange = $range.Duplicate\n$insertRange.Collapse(1)  # wdCollapseStart\n$insertRange.InsertBefore(\"los juegos de \")\n"}
